# The uploaded workbook removed two rows from the underlying source data:
#   - "Gasto executado"     (was in column F, row 9)
#   - "Data da liberação"   (was in column G, row 5)
# Excel's own row/column layout (A:G header columns) is untouched; only the
# label lists in columns F and G need their remaining entries shifted up by
# one to close the gap left by the removed rows, with the final cell in each
# column cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F ("Fonte FAPESP" list): remove "Gasto executado" at F9, shift the
# rest of the column up by one row.
for ($r = 9; $r -le 14; $r++) {
    $ws.Cells.Item($r, 6).Value2 = $ws.Cells.Item($r + 1, 6).Value2
}
$ws.Cells.Item(15, 6).Value2 = ""

# Column G ("Fonte FINEP" list): remove "Data da liberação" at G5, shift the
# rest of the column up by one row.
for ($r = 5; $r -le 14; $r++) {
    $ws.Cells.Item($r, 7).Value2 = $ws.Cells.Item($r + 1, 7).Value2
}
$ws.Cells.Item(15, 7).Value2 = ""

# Match the author's final selection in the saved file.
$ws.Range("F9").Select() | Out-Null
